$p = $ppt.ActivePresentation

# The deck's single Slide Master currently carries the "Integral" (Red
# Violet) theme color scheme. This edit swaps it back to the plain
# "Office Theme" color scheme (the other theme part already embedded in
# the package). We do this by rewriting the 12 theme colors through the
# ThemeColorScheme exposed on a slide - this is shared by every slide /
# the slide master, so one pass updates the whole design.

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# msoThemeColorIndex order: 1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
# 5-10 Accent1-6, 11 Hyperlink, 12 FollowedHyperlink.
# RGB is packed as 0x00BBGGRR (OLE COLORREF), so each target hex color
# "RRGGBB" becomes R | (G<<8) | (B<<16).
$tcs.Colors(1).RGB = 0          # Dark1    000000
$tcs.Colors(2).RGB = 16777215   # Light1   FFFFFF
$tcs.Colors(3).RGB = 6968388    # Dark2    44546A
$tcs.Colors(4).RGB = 15132391   # Light2   E7E6E6
$tcs.Colors(5).RGB = 13998939   # Accent1  5B9BD5
$tcs.Colors(6).RGB = 3243501    # Accent2  ED7D31
$tcs.Colors(7).RGB = 10855845   # Accent3  A5A5A5
$tcs.Colors(8).RGB = 49407      # Accent4  FFC000
$tcs.Colors(9).RGB = 12874308   # Accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # Accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # Hyperlink         0563C1
$tcs.Colors(12).RGB = 7491477   # FollowedHyperlink 954F72
